$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "zone" value shared by E2/E3 ("20 N" -> "32 T")
$ws.Range("E2").Value = "32 T"
$ws.Range("E3").Value = "32 T"

# Update installation / O&M related numeric figures
$ws.Range("C2").Value = 318989
$ws.Range("D2").Value = 4696615
$ws.Range("C3").Value = 296136
$ws.Range("D3").Value = 4641484

# Move the active selection
$ws.Range("D5").Select()
